$wb = $excel.ActiveWorkbook

$sheetNames = @("FOB", "DC", "RC")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("B2").Value = 4.402256559000539
    $ws.Range("F2").Value = 0.5209675996191799
    $ws.Range("J2").Value = 1470.008428551225
    $ws.Range("B3").Value = 3.687830186065481
    $ws.Range("F3").Value = 0.4310449375540863
    $ws.Range("J3").Value = 1626.13191320582
    $ws.Range("B4").Value = 3.819975580303203
    $ws.Range("F4").Value = 0.4487638864177364
    $ws.Range("J4").Value = 1589.66370027347
    $ws.Range("B5").Value = 3.716279683846471
    $ws.Range("F5").Value = 0.4307577569096745
    $ws.Range("J5").Value = 1633.480495801996
    $ws.Range("B6").Value = 3.661055697148117
    $ws.Range("F6").Value = 0.4072754945379497
    $ws.Range("J6").Value = 1714.777406735652
    $ws.Range("B7").Value = 3.585953026400534
    $ws.Range("F7").Value = 0.4004866486731615
    $ws.Range("J7").Value = 1725.866187685428
    $ws.Range("B8").Value = 3.507710792402266
    $ws.Range("F8").Value = 0.4063243739063441
    $ws.Range("J8").Value = 1682.410147361587
    $ws.Range("B9").Value = 3.865560412636638
    $ws.Range("F9").Value = 0.4392324838149131
    $ws.Range("J9").Value = 1633.821646415325
    $ws.Range("B10").Value = 3.516771457998645
    $ws.Range("F10").Value = 0.420695587734931
    $ws.Range("J10").Value = 1627.035318264224
    $ws.Range("B11").Value = 3.110579763551217
    $ws.Range("F11").Value = 0.3683028341410021
    $ws.Range("J11").Value = 1747.867369787412
    $ws.Range("B12").Value = 3.197935460579468
    $ws.Range("F12").Value = 0.3636140133157875
    $ws.Range("J12").Value = 1795.093604729728
    $ws.Range("B13").Value = 3.147449728872417
    $ws.Range("F13").Value = 0.3740954929883895
    $ws.Range("J13").Value = 1730.970990831482
    $ws.Range("B14").Value = 3.844423804335051
    $ws.Range("F14").Value = 0.4482410056471485
    $ws.Range("J14").Value = 1596.602884262822
    $ws.Range("B15").Value = 3.048251498678792
    $ws.Range("F15").Value = 0.3816450771064601
    $ws.Range("J15").Value = 1669.777473168846
    $ws.Range("B16").Value = 3.868943027753081
    $ws.Range("F16").Value = 0.4144396832237783
    $ws.Range("J16").Value = 1732.318321579272
    $ws.Range("B17").Value = 2.692550679356606
    $ws.Range("F17").Value = 0.3203722225268264
    $ws.Range("J17").Value = 1869.476255797456
    $ws.Range("B18").Value = 3.147701475878357
    $ws.Range("F18").Value = 0.3779425077894463
    $ws.Range("J18").Value = 1713.420240378927
    $ws.Range("B19").Value = 2.772632298534591
    $ws.Range("F19").Value = 0.3475241814562106
    $ws.Range("J19").Value = 1748.855690119198
    $ws.Range("B20").Value = 3.326774862942547
    $ws.Range("F20").Value = 0.4180135463247534
    $ws.Range("J20").Value = 1592.62751366822
    $ws.Range("B21").Value = 3.269341196041328
    $ws.Range("F21").Value = 0.397070607318797
    $ws.Range("J21").Value = 1662.092721981384
    $ws.Range("B22").Value = 2.870472262021966
    $ws.Range("F22").Value = 0.3517771549135828
    $ws.Range("J22").Value = 1757.931343644857
    $ws.Range("B23").Value = 3.891625364726056
    $ws.Range("F23").Value = 0.4266201064407552
    $ws.Range("J23").Value = 1687.784745792284
    $ws.Range("B24").Value = 2.580189589162795
    $ws.Range("F24").Value = 0.3238330606615358
    $ws.Range("J24").Value = 1810.495664627731
    $ws.Range("B25").Value = 3.170239468367399
    $ws.Range("F25").Value = 0.3616976584486941
    $ws.Range("J25").Value = 1796.772952858149
    $ws.Range("B26").Value = 2.654357158473029
    $ws.Range("F26").Value = 0.3069015387608979
    $ws.Range("J26").Value = 1937.641725741084
}
